$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 100018.7
$ws.Range("I11").Value = 100018.7
$ws.Range("K11").Value = 100018.7
$ws.Range("M11").Value = -99878.7
$ws.Range("H33").Value = 371680.84
$ws.Range("I33").Value = 472.7857
$ws.Range("J33").Value = 1114097.0
$ws.Range("K33").Value = 472.7857
$ws.Range("L33").Value = 1114097.0
$ws.Range("M33").Value = -243.7857
$ws.Range("N33").Value = -1114555.0
$ws.Range("H98").Value = 732.1429
$ws.Range("I98").Value = 525.5
$ws.Range("J98").Value = 1007.6667
$ws.Range("K98").Value = 525.5
$ws.Range("L98").Value = 1007.6667
$ws.Range("M98").Value = 972.5
$ws.Range("N98").Value = -4003.6667
$ws.Range("H122").Value = 732.1429
$ws.Range("I122").Value = 525.5
$ws.Range("J122").Value = 1007.6667
$ws.Range("K122").Value = 1576.5
$ws.Range("L122").Value = 3023.0001
$ws.Range("M122").Value = 873.5
$ws.Range("N122").Value = -7923.0001
$ws.Range("H129").Value = 2532.1865
$ws.Range("I129").Value = 11560.667
$ws.Range("J129").Value = 907.06
$ws.Range("K129").Value = 34682.001
$ws.Range("L129").Value = 2721.18
$ws.Range("M129").Value = -29682.001
$ws.Range("N129").Value = -12721.18
$ws.Range("H138").Value = 2991.2222
$ws.Range("I138").Value = 1638.069
$ws.Range("J138").Value = 3551.8142
$ws.Range("K138").Value = 4914.207
$ws.Range("L138").Value = 10655.4426
$ws.Range("M138").Value = 225.7929999999997
$ws.Range("N138").Value = -20935.4426

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29553.158
$ws.Range("I32").Value = 10666.845
$ws.Range("K32").Value = 10666.845
$ws.Range("M32").Value = -10379.845
$ws.Range("H61").Value = 1536.5333
$ws.Range("I61").Value = 966.0
$ws.Range("J61").Value = 2107.0667
$ws.Range("K61").Value = 966.0
$ws.Range("L61").Value = 2107.0667
$ws.Range("M61").Value = -754.0
$ws.Range("N61").Value = -2531.0667
$ws.Range("H74").Value = 1624.6
$ws.Range("I74").Value = 1936.7
$ws.Range("J74").Value = 1312.5
$ws.Range("K74").Value = 1936.7
$ws.Range("L74").Value = 1312.5
$ws.Range("M74").Value = -1062.7
$ws.Range("N74").Value = -3060.5
$ws.Range("H77").Value = 1624.6
$ws.Range("I77").Value = 1936.7
$ws.Range("J77").Value = 1312.5
$ws.Range("K77").Value = 9683.5
$ws.Range("L77").Value = 6562.5
$ws.Range("M77").Value = -5315.5
$ws.Range("N77").Value = -15298.5
$ws.Range("H80").Value = 25252.445
$ws.Range("J80").Value = 25252.445
$ws.Range("L80").Value = 25252.445
$ws.Range("N80").Value = -27248.445
$ws.Range("H83").Value = 25252.445
$ws.Range("J83").Value = 25252.445
$ws.Range("L83").Value = 75757.33499999999
$ws.Range("N83").Value = -85741.33499999999
$ws.Range("H122").Value = 2976.0667
$ws.Range("I122").Value = 2088.1
$ws.Range("K122").Value = 6264.299999999999
$ws.Range("M122").Value = -3814.299999999999
$ws.Range("H131").Value = 25412.785
$ws.Range("J131").Value = 25412.785
$ws.Range("L131").Value = 25412.785
$ws.Range("N131").Value = -35492.785
$ws.Range("H132").Value = 14787.761
$ws.Range("I132").Value = 17195.764
$ws.Range("J132").Value = 3349.75
$ws.Range("K132").Value = 51587.292
$ws.Range("L132").Value = 10049.25
$ws.Range("M132").Value = -49057.292
$ws.Range("N132").Value = -15109.25
$ws.Range("H136").Value = 1536.5333
$ws.Range("I136").Value = 966.0
$ws.Range("J136").Value = 2107.0667
$ws.Range("K136").Value = 2898.0
$ws.Range("L136").Value = 6321.2001
$ws.Range("M136").Value = -348.0
$ws.Range("N136").Value = -11421.2001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2415.8635
$ws.Range("I134").Value = 2242.2593
$ws.Range("J134").Value = 3197.0833
$ws.Range("K134").Value = 6726.777900000001
$ws.Range("L134").Value = 9591.249899999999
$ws.Range("M134").Value = -4191.777900000001
$ws.Range("N134").Value = -14661.2499

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 6000000.0
$ws.Range("I4").Value = 10000000.0
$ws.Range("K4").Value = 10000000.0
$ws.Range("M4").Value = -9999888.0
$ws.Range("H58").Value = 2429.7646
$ws.Range("I58").Value = 2342.6667
$ws.Range("J58").Value = 2638.8
$ws.Range("K58").Value = 2342.6667
$ws.Range("L58").Value = 2638.8
$ws.Range("M58").Value = -2139.6667
$ws.Range("N58").Value = -3044.8
$ws.Range("H94").Value = 1284.7273
$ws.Range("J94").Value = 1309.5625
$ws.Range("L94").Value = 1309.5625
$ws.Range("N94").Value = -2211.5625
$ws.Range("H99").Value = 2069.5925
$ws.Range("I99").Value = 1529.0
$ws.Range("J99").Value = 2387.5881
$ws.Range("K99").Value = 1529.0
$ws.Range("L99").Value = 2387.5881
$ws.Range("M99").Value = -31.0
$ws.Range("N99").Value = -5383.5881
$ws.Range("H107").Value = 784.5
$ws.Range("I107").Value = 741.4
$ws.Range("K107").Value = 741.4
$ws.Range("M107").Value = 1178.6
$ws.Range("H122").Value = 970.5
$ws.Range("I122").Value = 775.0
$ws.Range("J122").Value = 1035.6666
$ws.Range("K122").Value = 2325.0
$ws.Range("L122").Value = 3106.9998
$ws.Range("M122").Value = 125.0
$ws.Range("N122").Value = -8006.9998
$ws.Range("H126").Value = 2069.5925
$ws.Range("I126").Value = 1529.0
$ws.Range("J126").Value = 2387.5881
$ws.Range("K126").Value = 4587.0
$ws.Range("L126").Value = 7162.7643
$ws.Range("M126").Value = -2117.0
$ws.Range("N126").Value = -12102.7643
$ws.Range("H132").Value = 1934.0
$ws.Range("I132").Value = 1852.591
$ws.Range("J132").Value = 2209.5386
$ws.Range("K132").Value = 5557.772999999999
$ws.Range("L132").Value = 6628.6158
$ws.Range("M132").Value = -3027.772999999999
$ws.Range("N132").Value = -11688.6158
$ws.Range("H134").Value = 2089.111
$ws.Range("I134").Value = 1049.75
$ws.Range("J134").Value = 2920.6
$ws.Range("K134").Value = 3149.25
$ws.Range("L134").Value = 8761.8
$ws.Range("M134").Value = -614.25
$ws.Range("N134").Value = -13831.8
$ws.Range("H136").Value = 2429.7646
$ws.Range("I136").Value = 2342.6667
$ws.Range("J136").Value = 2638.8
$ws.Range("K136").Value = 7028.000100000001
$ws.Range("L136").Value = 7916.400000000001
$ws.Range("M136").Value = -4478.000100000001
$ws.Range("N136").Value = -13016.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 34.863636
$ws.Range("I12").Value = 11.6
$ws.Range("J12").Value = 41.705883
$ws.Range("K12").Value = 34.8
$ws.Range("L12").Value = 125.117649
$ws.Range("M12").Value = 138.2
$ws.Range("N12").Value = -471.117649
$ws.Range("H37").Value = 647505.6
$ws.Range("J37").Value = 647505.6
$ws.Range("L37").Value = 1942516.8
$ws.Range("N37").Value = -1942740.8
$ws.Range("H92").Value = 3003.0
$ws.Range("I92").Value = 0.0
$ws.Range("K92").Value = 0.0
$ws.Range("M92").ClearContents()
$ws.Range("H98").Value = 91928.73
$ws.Range("J98").Value = 101121.3
$ws.Range("L98").Value = 303363.9
$ws.Range("N98").Value = -306359.9
$ws.Range("H109").Value = 3342.6858
$ws.Range("I109").Value = 2725.6667
$ws.Range("J109").Value = 3470.3447
$ws.Range("K109").Value = 8177.000100000001
$ws.Range("L109").Value = 10411.0341
$ws.Range("M109").Value = -7137.000100000001
$ws.Range("N109").Value = -12491.0341
$ws.Range("H131").Value = 810.83
$ws.Range("J131").Value = 870.2651
$ws.Range("L131").Value = 2610.7953
$ws.Range("N131").Value = -12690.7953

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 5291429.0
$ws.Range("I11").Value = 5291429.0
$ws.Range("J11").Value = 0.0
$ws.Range("K11").Value = 5291429.0
$ws.Range("L11").Value = 0.0
$ws.Range("M11").Value = -5291290.0
$ws.Range("N11").ClearContents()
$ws.Range("H70").Value = 173908.92
$ws.Range("I70").Value = 256299.75
$ws.Range("K70").Value = 256299.75
$ws.Range("M70").Value = -256029.75
$ws.Range("H73").Value = 173908.92
$ws.Range("I73").Value = 256299.75
$ws.Range("K73").Value = 256299.75
$ws.Range("M73").Value = -255363.75
$ws.Range("H96").Value = 0.0
$ws.Range("J96").Value = 0.0
$ws.Range("L96").Value = 0.0
$ws.Range("N96").ClearContents()
$ws.Range("H102").Value = 302344.4
$ws.Range("I102").Value = 2183.6924
$ws.Range("J102").Value = 859785.7
$ws.Range("K102").Value = 2183.6924
$ws.Range("L102").Value = 859785.7
$ws.Range("M102").Value = -561.6923999999999
$ws.Range("N102").Value = -863029.7
$ws.Range("H132").Value = 3116.1428
$ws.Range("I132").Value = 2345.1875
$ws.Range("J132").Value = 4144.0835
$ws.Range("K132").Value = 7035.5625
$ws.Range("L132").Value = 12432.2505
$ws.Range("M132").Value = -4505.5625
$ws.Range("N132").Value = -17492.2505

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4178.7915
$ws.Range("I132").Value = 5891.636
$ws.Range("J132").Value = 2729.4614
$ws.Range("K132").Value = 17674.908
$ws.Range("L132").Value = 8188.3842
$ws.Range("M132").Value = -15144.908
$ws.Range("N132").Value = -13248.3842

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 52415.0
$ws.Range("I26").Value = 4830.5
$ws.Range("J26").Value = 99999.5
$ws.Range("K26").Value = 4830.5
$ws.Range("L26").Value = 99999.5
$ws.Range("M26").Value = -4537.5
$ws.Range("N26").Value = -100585.5
$ws.Range("H81").Value = 168895.75
$ws.Range("I81").Value = 112194.445
$ws.Range("J81").Value = 338999.66
$ws.Range("K81").Value = 224388.89
$ws.Range("L81").Value = 677999.32
$ws.Range("M81").Value = -223327.89
$ws.Range("N81").Value = -680121.32
$ws.Range("H84").Value = 168895.75
$ws.Range("I84").Value = 112194.445
$ws.Range("J84").Value = 338999.66
$ws.Range("K84").Value = 1121944.45
$ws.Range("L84").Value = 3389996.6
$ws.Range("M84").Value = -1116640.45
$ws.Range("N84").Value = -3400604.6
$ws.Range("H122").Value = 1787.5625
$ws.Range("I122").Value = 1828.6428
$ws.Range("J122").Value = 1500.0
$ws.Range("K122").Value = 5485.928400000001
$ws.Range("L122").Value = 4500.0
$ws.Range("M122").Value = -3035.928400000001
$ws.Range("N122").Value = -9400.0
$ws.Range("H132").Value = 3208.957
$ws.Range("I132").Value = 1637.5283
$ws.Range("K132").Value = 4912.5849
$ws.Range("M132").Value = -2382.5849
